# repull data, push all data, mean calculation
# Update the "dSF" (column F) values for a set of rows to reflect repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    4  = -3
    6  = 2
    8  = 2
    9  = 2
    11 = -2
    15 = 0
    22 = 3
    24 = -4
    31 = 1
    35 = -3
    38 = 3
    42 = -7
    43 = -5
    44 = 0
    46 = -1
    49 = 0
    50 = 0
    60 = 0
    61 = 1
    62 = -1
    65 = 1
    70 = 0
    76 = 2
    84 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
